$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 375
$ws.Range("J8").Value = 1325.6
$ws.Range("L8").Value = 3976.8
$ws.Range("N8").Value = -4254.799999999999
$ws.Range("H11").Value = 41.285713
$ws.Range("I11").Value = 41.285713
$ws.Range("K11").Value = 41.285713
$ws.Range("M11").Value = 98.714287
$ws.Range("H17").Value = 2074090
$ws.Range("I17").Value = 705
$ws.Range("K17").Value = 2115
$ws.Range("M17").Value = -1947
$ws.Range("H19").Value = 23935.834
$ws.Range("I19").Value = 534
$ws.Range("J19").Value = 47337.668
$ws.Range("K19").Value = 534
$ws.Range("L19").Value = 47337.668
$ws.Range("M19").Value = -359
$ws.Range("N19").Value = -47687.668
$ws.Range("H28").Value = 501.6129
$ws.Range("I28").Value = 501.69232
$ws.Range("K28").Value = 501.69232
$ws.Range("M28").Value = -16.69232
$ws.Range("H55").Value = 398.5
$ws.Range("I55").Value = 398.5
$ws.Range("K55").Value = 398.5
$ws.Range("M55").Value = -184.5
$ws.Range("H100").Value = 4725.4
$ws.Range("I100").Value = 2066.25
$ws.Range("J100").Value = 6498.1665
$ws.Range("K100").Value = 2066.25
$ws.Range("L100").Value = 6498.1665
$ws.Range("M100").Value = -1525.25
$ws.Range("N100").Value = -7580.1665
$ws.Range("H107").Value = 681.1667
$ws.Range("I107").Value = 419.63635
$ws.Range("K107").Value = 419.63635
$ws.Range("M107").Value = 1500.36365
$ws.Range("H111").Value = 675
$ws.Range("I111").Value = 520.1818
$ws.Range("J111").Value = 1100.75
$ws.Range("K111").Value = 1560.5454
$ws.Range("L111").Value = 3302.25
$ws.Range("M111").Value = 1506.4546
$ws.Range("N111").Value = -9436.25
$ws.Range("H112").Value = 2210.2727
$ws.Range("I112").Value = 450
$ws.Range("J112").Value = 2386.3
$ws.Range("K112").Value = 1350
$ws.Range("L112").Value = 7158.900000000001
$ws.Range("M112").Value = -242
$ws.Range("N112").Value = -9374.900000000001
$ws.Range("H123").Value = 113555.8
$ws.Range("J123").Value = 113555.8
$ws.Range("L123").Value = 113555.8
$ws.Range("N123").Value = -123355.8
$ws.Range("H125").Value = 6298.3335
$ws.Range("J125").Value = 7379.8
$ws.Range("L125").Value = 66418.2
$ws.Range("N125").Value = -71338.2
$ws.Range("H137").Value = 12692.156
$ws.Range("I137").Value = 3362.4
$ws.Range("J137").Value = 20924.295
$ws.Range("K137").Value = 10087.2
$ws.Range("L137").Value = 62772.88499999999
$ws.Range("M137").Value = -7537.200000000001
$ws.Range("N137").Value = -67872.88499999999
$ws.Range("H138").Value = 3371.4285
$ws.Range("J138").Value = 2519.1538
$ws.Range("L138").Value = 7557.4614
$ws.Range("N138").Value = -17837.4614
$ws.Range("H141").Value = 3850

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21861.533
$ws.Range("I61").Value = 10473.75
$ws.Range("J61").Value = 26002.545
$ws.Range("K61").Value = 10473.75
$ws.Range("L61").Value = 26002.545
$ws.Range("M61").Value = -10261.75
$ws.Range("N61").Value = -26426.545
$ws.Range("H63").Value = 5587.2
$ws.Range("J63").Value = 3499.5
$ws.Range("L63").Value = 3499.5
$ws.Range("N63").Value = -4871.5
$ws.Range("H66").Value = 5587.2
$ws.Range("J66").Value = 3499.5
$ws.Range("L66").Value = 17497.5
$ws.Range("N66").Value = -24361.5
$ws.Range("H74").Value = 18665.2
$ws.Range("I74").Value = 2562.8125
$ws.Range("J74").Value = 47291.668
$ws.Range("K74").Value = 2562.8125
$ws.Range("L74").Value = 47291.668
$ws.Range("M74").Value = -1688.8125
$ws.Range("N74").Value = -49039.668
$ws.Range("H77").Value = 18665.2
$ws.Range("I77").Value = 2562.8125
$ws.Range("J77").Value = 47291.668
$ws.Range("K77").Value = 12814.0625
$ws.Range("L77").Value = 236458.34
$ws.Range("M77").Value = -8446.0625
$ws.Range("N77").Value = -245194.34
$ws.Range("H88").Value = 2346
$ws.Range("I88").Value = 1250
$ws.Range("K88").Value = 1250
$ws.Range("M88").Value = -844
$ws.Range("H91").Value = 2346
$ws.Range("I91").Value = 1250
$ws.Range("K91").Value = 1250
$ws.Range("M91").Value = 154
$ws.Range("H97").Value = 2239.238
$ws.Range("I97").Value = 1966.1765
$ws.Range("K97").Value = 1966.1765
$ws.Range("M97").Value = -1470.1765
$ws.Range("H102").Value = 16506.732
$ws.Range("I102").Value = 3400.1428
$ws.Range("K102").Value = 3400.1428
$ws.Range("M102").Value = -1778.1428
$ws.Range("H122").Value = 6485.2856
$ws.Range("I122").Value = 4465.6665
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 13396.9995
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -10946.9995
$ws.Range("N122").Value = -28900
$ws.Range("H130").Value = 77487.625
$ws.Range("J130").Value = 77487.625
$ws.Range("L130").Value = 77487.625
$ws.Range("N130").Value = -87527.625
$ws.Range("H132").Value = 9119010
$ws.Range("I132").Value = 18727
$ws.Range("K132").Value = 56181
$ws.Range("M132").Value = -53651
$ws.Range("H136").Value = 21861.533
$ws.Range("I136").Value = 10473.75
$ws.Range("J136").Value = 26002.545
$ws.Range("K136").Value = 31421.25
$ws.Range("L136").Value = 78007.63499999999
$ws.Range("M136").Value = -28871.25
$ws.Range("N136").Value = -83107.63499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2069.1
$ws.Range("J80").Value = 2221.2222
$ws.Range("L80").Value = 2221.2222
$ws.Range("N80").Value = -4217.2222
$ws.Range("H83").Value = 2069.1
$ws.Range("J83").Value = 2221.2222
$ws.Range("L83").Value = 11106.111
$ws.Range("N83").Value = -21090.111
$ws.Range("H86").Value = 1241
$ws.Range("I86").Value = 1342.8462
$ws.Range("J86").Value = 976.2
$ws.Range("K86").Value = 1342.8462
$ws.Range("L86").Value = 976.2
$ws.Range("M86").Value = -219.8462
$ws.Range("N86").Value = -3222.2
$ws.Range("H89").Value = 1241
$ws.Range("I89").Value = 1342.8462
$ws.Range("J89").Value = 976.2
$ws.Range("K89").Value = 6714.231
$ws.Range("L89").Value = 4881
$ws.Range("M89").Value = -1098.231
$ws.Range("N89").Value = -16113
$ws.Range("H94").Value = 4549.7085
$ws.Range("I94").Value = 3806.2058
$ws.Range("K94").Value = 3806.2058
$ws.Range("M94").Value = -3355.2058
$ws.Range("H99").Value = 14056.6875
$ws.Range("J99").Value = 31314.143
$ws.Range("L99").Value = 31314.143
$ws.Range("N99").Value = -34310.143
$ws.Range("H105").Value = 1462.6765
$ws.Range("I105").Value = 1164.96
$ws.Range("J105").Value = 2289.6667
$ws.Range("K105").Value = 1164.96
$ws.Range("L105").Value = 2289.6667
$ws.Range("M105").Value = 582.04
$ws.Range("N105").Value = -5783.6667
$ws.Range("H107").Value = 2539.652
$ws.Range("I107").Value = 2291.7334
$ws.Range("K107").Value = 2291.7334
$ws.Range("M107").Value = -371.7334000000001
$ws.Range("H134").Value = 10439.081
$ws.Range("I134").Value = 5073.4136
$ws.Range("J134").Value = 29889.625
$ws.Range("K134").Value = 15220.2408
$ws.Range("L134").Value = 89668.875
$ws.Range("M134").Value = -12685.2408
$ws.Range("N134").Value = -94738.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1529.6428
$ws.Range("I22").Value = 597.5
$ws.Range("J22").Value = 1685
$ws.Range("K22").Value = 597.5
$ws.Range("L22").Value = 1685
$ws.Range("M22").Value = -247.5
$ws.Range("N22").Value = -2385
$ws.Range("H31").Value = 28943.479
$ws.Range("I31").Value = 15975.714
$ws.Range("J31").Value = 34616.875
$ws.Range("K31").Value = 15975.714
$ws.Range("L31").Value = 34616.875
$ws.Range("M31").Value = -15680.714
$ws.Range("N31").Value = -35206.875
$ws.Range("H34").Value = 28943.479
$ws.Range("I34").Value = 15975.714
$ws.Range("J34").Value = 34616.875
$ws.Range("K34").Value = 15975.714
$ws.Range("L34").Value = 34616.875
$ws.Range("M34").Value = -15773.714
$ws.Range("N34").Value = -35020.875
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null
$ws.Range("H58").Value = 16630.31
$ws.Range("J58").Value = 17062.791
$ws.Range("L58").Value = 17062.791
$ws.Range("N58").Value = -17468.791
$ws.Range("H70").Value = 17000
$ws.Range("J70").Value = 17000
$ws.Range("L70").Value = 17000
$ws.Range("N70").Value = -17630
$ws.Range("H73").Value = 17000
$ws.Range("J73").Value = 17000
$ws.Range("L73").Value = 17000
$ws.Range("N73").Value = -19184
$ws.Range("H99").Value = 6537
$ws.Range("I99").Value = 3441.5
$ws.Range("J99").Value = 8788.272000000001
$ws.Range("K99").Value = 3441.5
$ws.Range("L99").Value = 8788.272000000001
$ws.Range("M99").Value = -1943.5
$ws.Range("N99").Value = -11784.272
$ws.Range("H105").Value = 13353.5
$ws.Range("I105").Value = 25707
$ws.Range("K105").Value = 25707
$ws.Range("M105").Value = -23960
$ws.Range("H107").Value = 1683.8
$ws.Range("I107").Value = 1128.5
$ws.Range("J107").Value = 3288
$ws.Range("K107").Value = 1128.5
$ws.Range("L107").Value = 3288
$ws.Range("M107").Value = 791.5
$ws.Range("N107").Value = -7128
$ws.Range("H122").Value = 3502.6
$ws.Range("I122").Value = 2834.3333
$ws.Range("K122").Value = 8502.999899999999
$ws.Range("M122").Value = -6052.999899999999
$ws.Range("H126").Value = 6537
$ws.Range("I126").Value = 3441.5
$ws.Range("J126").Value = 8788.272000000001
$ws.Range("K126").Value = 10324.5
$ws.Range("L126").Value = 26364.816
$ws.Range("M126").Value = -7854.5
$ws.Range("N126").Value = -31304.816
$ws.Range("H133").Value = 80525.91
$ws.Range("J133").Value = 80525.91
$ws.Range("L133").Value = 80525.91
$ws.Range("N133").Value = -85585.91
$ws.Range("H134").Value = 41675708
$ws.Range("I134").Value = 2294.4285
$ws.Range("J134").Value = 58835350
$ws.Range("K134").Value = 6883.2855
$ws.Range("L134").Value = 176506050
$ws.Range("M134").Value = -4348.2855
$ws.Range("N134").Value = -176511120
$ws.Range("H136").Value = 16630.31
$ws.Range("J136").Value = 17062.791
$ws.Range("L136").Value = 51188.37300000001
$ws.Range("N136").Value = -56288.37300000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2258.6667
$ws.Range("I11").Value = 2679.75
$ws.Range("J11").Value = 1416.5
$ws.Range("K11").Value = 8039.25
$ws.Range("L11").Value = 4249.5
$ws.Range("M11").Value = -7899.25
$ws.Range("N11").Value = -4529.5
$ws.Range("H12").Value = 325
$ws.Range("I12").Value = 349.5
$ws.Range("J12").Value = 300.5
$ws.Range("K12").Value = 1048.5
$ws.Range("L12").Value = 901.5
$ws.Range("M12").Value = -875.5
$ws.Range("N12").Value = -1247.5
$ws.Range("H16").Value = 5.6666665
$ws.Range("I16").Value = 1
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 170
$ws.Range("H17").Value = 238.16667
$ws.Range("I17").Value = 256.875
$ws.Range("J17").Value = 223.2
$ws.Range("K17").Value = 770.625
$ws.Range("L17").Value = 669.5999999999999
$ws.Range("M17").Value = -601.625
$ws.Range("N17").Value = -1007.6
$ws.Range("H19").Value = 167.33333
$ws.Range("J19").Value = 51
$ws.Range("L19").Value = 153
$ws.Range("N19").Value = -501
$ws.Range("H20").Value = 1593.2667
$ws.Range("J20").Value = 1637.5
$ws.Range("L20").Value = 4912.5
$ws.Range("N20").Value = -5366.5
$ws.Range("H38").Value = 113.5
$ws.Range("I38").Value = 16.75
$ws.Range("J38").Value = 178
$ws.Range("K38").Value = 50.25
$ws.Range("L38").Value = 534
$ws.Range("M38").Value = 296.75
$ws.Range("N38").Value = -1228
$ws.Range("H39").Value = 3482.1667
$ws.Range("J39").Value = 3497
$ws.Range("L39").Value = 10491
$ws.Range("N39").Value = -11079
$ws.Range("H45").Value = 766.5
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -968
$ws.Range("H46").Value = 315.66666
$ws.Range("I46").Value = 466.66666
$ws.Range("J46").Value = 164.66667
$ws.Range("K46").Value = 1399.99998
$ws.Range("L46").Value = 494.00001
$ws.Range("M46").Value = -1308.99998
$ws.Range("N46").Value = -676.00001
$ws.Range("H50").Value = 2194660.2
$ws.Range("I50").Value = 3201.5
$ws.Range("J50").Value = 4386119
$ws.Range("K50").Value = 9604.5
$ws.Range("L50").Value = 13158357
$ws.Range("M50").Value = -9123.5
$ws.Range("N50").Value = -13159319
$ws.Range("H53").Value = 2194660.2
$ws.Range("I53").Value = 3201.5
$ws.Range("J53").Value = 4386119
$ws.Range("K53").Value = 9604.5
$ws.Range("L53").Value = 13158357
$ws.Range("M53").Value = -9123.5
$ws.Range("N53").Value = -13159319
$ws.Range("H57").Value = 2000
$ws.Range("J57").Value = 1000
$ws.Range("L57").Value = 3000
$ws.Range("N57").Value = -4118
$ws.Range("H59").Value = 101652.2
$ws.Range("I59").Value = 1198.5
$ws.Range("J59").Value = 168621.33
$ws.Range("K59").Value = 3595.5
$ws.Range("L59").Value = 505863.99
$ws.Range("M59").Value = -3055.5
$ws.Range("N59").Value = -506943.99
$ws.Range("H109").Value = 2224199
$ws.Range("I109").Value = 1822.1818
$ws.Range("J109").Value = 8335735.5
$ws.Range("K109").Value = 5466.5454
$ws.Range("L109").Value = 25007206.5
$ws.Range("M109").Value = -4426.5454
$ws.Range("N109").Value = -25009286.5
$ws.Range("H117").Value = 1978.875
$ws.Range("J117").Value = 2218.7144
$ws.Range("L117").Value = 6656.1432
$ws.Range("N117").Value = -13540.1432
$ws.Range("H120").Value = 25168.334
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = $null
$ws.Range("H122").Value = 8609066
$ws.Range("I122").Value = 11680082
$ws.Range("K122").Value = 105120738
$ws.Range("M122").Value = -105118288
$ws.Range("H131").Value = 1479.51
$ws.Range("I131").Value = 1014.5
$ws.Range("J131").Value = 1489
$ws.Range("K131").Value = 3043.5
$ws.Range("L131").Value = 4467
$ws.Range("M131").Value = 1996.5
$ws.Range("N131").Value = -14547
$ws.Range("H137").Value = 1518.1
$ws.Range("J137").Value = 1783.1666
$ws.Range("L137").Value = 5349.4998
$ws.Range("N137").Value = -15549.4998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 62
$ws.Range("H70").Value = 11961.619
$ws.Range("I70").Value = 17850.75
$ws.Range("J70").Value = 8337.538
$ws.Range("K70").Value = 17850.75
$ws.Range("L70").Value = 8337.538
$ws.Range("M70").Value = -17580.75
$ws.Range("N70").Value = -8877.538
$ws.Range("H73").Value = 11961.619
$ws.Range("I73").Value = 17850.75
$ws.Range("J73").Value = 8337.538
$ws.Range("K73").Value = 17850.75
$ws.Range("L73").Value = 8337.538
$ws.Range("M73").Value = -16914.75
$ws.Range("N73").Value = -10209.538
$ws.Range("H80").Value = 8091.963
$ws.Range("I80").Value = 4563.6313
$ws.Range("J80").Value = 16471.75
$ws.Range("K80").Value = 4563.6313
$ws.Range("L80").Value = 16471.75
$ws.Range("M80").Value = -3565.6313
$ws.Range("N80").Value = -18467.75
$ws.Range("H83").Value = 8091.963
$ws.Range("I83").Value = 4563.6313
$ws.Range("J83").Value = 16471.75
$ws.Range("K83").Value = 22818.1565
$ws.Range("L83").Value = 82358.75
$ws.Range("M83").Value = -17826.1565
$ws.Range("N83").Value = -92342.75
$ws.Range("H107").Value = 892.55554
$ws.Range("I107").Value = 441.83334
$ws.Range("J107").Value = 1117.9166
$ws.Range("K107").Value = 441.83334
$ws.Range("L107").Value = 1117.9166
$ws.Range("M107").Value = 1478.16666
$ws.Range("N107").Value = -4957.9166
$ws.Range("H113").Value = 51737.277
$ws.Range("I113").Value = 71067.53999999999
$ws.Range("J113").Value = 1478.6
$ws.Range("K113").Value = 71067.53999999999
$ws.Range("L113").Value = 1478.6
$ws.Range("M113").Value = -68897.53999999999
$ws.Range("N113").Value = -5818.6
$ws.Range("H122").Value = 1404.75
$ws.Range("I122").Value = 1265.4286
$ws.Range("J122").Value = 1599.8
$ws.Range("K122").Value = 3796.2858
$ws.Range("L122").Value = 4799.4
$ws.Range("M122").Value = -1346.2858
$ws.Range("N122").Value = -9699.4
$ws.Range("H132").Value = 12527.421
$ws.Range("I132").Value = 8780.333000000001
$ws.Range("J132").Value = 18951
$ws.Range("K132").Value = 26340.999
$ws.Range("L132").Value = 56853
$ws.Range("M132").Value = -23810.999
$ws.Range("N132").Value = -61913

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").Value = $null
$ws.Range("H7").Value = 6994.4
$ws.Range("I7").Value = 3982.2
$ws.Range("K7").Value = 3982.2
$ws.Range("M7").Value = -3870.2
$ws.Range("H61").Value = 3240.7097
$ws.Range("I61").Value = 2118.1667
$ws.Range("J61").Value = 4795
$ws.Range("K61").Value = 2118.1667
$ws.Range("L61").Value = 4795
$ws.Range("M61").Value = -1916.1667
$ws.Range("N61").Value = -5199
$ws.Range("H82").Value = 5976.222
$ws.Range("I82").Value = 5414.75
$ws.Range("J82").Value = 7099.1665
$ws.Range("K82").Value = 5414.75
$ws.Range("L82").Value = 7099.1665
$ws.Range("M82").Value = -5053.75
$ws.Range("N82").Value = -7821.1665
$ws.Range("H85").Value = 5976.222
$ws.Range("I85").Value = 5414.75
$ws.Range("J85").Value = 7099.1665
$ws.Range("K85").Value = 5414.75
$ws.Range("L85").Value = 7099.1665
$ws.Range("M85").Value = -4166.75
$ws.Range("N85").Value = -9595.166499999999
$ws.Range("H113").Value = 3240.7097
$ws.Range("I113").Value = 2118.1667
$ws.Range("J113").Value = 4795
$ws.Range("K113").Value = 2118.1667
$ws.Range("L113").Value = 4795
$ws.Range("M113").Value = 51.83329999999978
$ws.Range("N113").Value = -9135
$ws.Range("H122").Value = 6178.615
$ws.Range("I122").Value = 5722
$ws.Range("J122").Value = 6570
$ws.Range("K122").Value = 17166
$ws.Range("L122").Value = 19710
$ws.Range("M122").Value = -14716
$ws.Range("N122").Value = -24610
$ws.Range("H126").Value = 6994.4
$ws.Range("I126").Value = 3982.2
$ws.Range("K126").Value = 11946.6
$ws.Range("M126").Value = -9476.599999999999
$ws.Range("H132").Value = 1220110.8
$ws.Range("I132").Value = 4526.875
$ws.Range("K132").Value = 13580.625
$ws.Range("M132").Value = -11050.625
$ws.Range("H136").Value = 10611.83
$ws.Range("I136").Value = 11147.625
$ws.Range("J136").Value = 10052.739
$ws.Range("K136").Value = 33442.875
$ws.Range("L136").Value = 30158.217
$ws.Range("M136").Value = -30892.875
$ws.Range("N136").Value = -35258.217

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").Value = $null
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").Value = $null
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("H107").Value = 1010.2414
$ws.Range("I107").Value = 937
$ws.Range("J107").Value = 1202.5
$ws.Range("K107").Value = 2811
$ws.Range("L107").Value = 3607.5
$ws.Range("M107").Value = -891
$ws.Range("N107").Value = -7447.5
$ws.Range("H122").Value = 4270.64
$ws.Range("I122").Value = 2834.3635
$ws.Range("J122").Value = 5399.143
$ws.Range("K122").Value = 8503.0905
$ws.Range("L122").Value = 16197.429
$ws.Range("M122").Value = -6053.0905
$ws.Range("N122").Value = -21097.429
$ws.Range("H126").Value = 16208.913
$ws.Range("I126").Value = 18440.6
$ws.Range("J126").Value = 1331
$ws.Range("K126").Value = 55321.8
$ws.Range("L126").Value = 3993
$ws.Range("M126").Value = -52851.8
$ws.Range("N126").Value = -8933
$ws.Range("H136").Value = 13065.667
$ws.Range("I136").Value = 2197.4167
$ws.Range("J136").Value = 27556.666
$ws.Range("K136").Value = 6592.250100000001
$ws.Range("L136").Value = 82669.99800000001
$ws.Range("M136").Value = -4042.250100000001
$ws.Range("N136").Value = -87769.99800000001
$ws.Range("H138").Value = 142497.25
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 142497.25
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = -152777.25

Write-Host "Applied all cell updates"